$wb = $excel.ActiveWorkbook

# ---- Sheet: PIR (A1:F277 -> A1:F291) ----
$ws = $wb.Worksheets.Item("PIR")
$rows = New-Object System.Collections.ArrayList
$rows.Add(@("2026-01-30", "17:30:12", "17:00", "Bathroom", "No Motion", "Inactive")) | Out-Null
$rows.Add(@("2026-01-30", "17:30:13", "17:00", "Bathroom", "No Motion", "Inactive")) | Out-Null
$rows.Add(@("2026-01-30", "17:30:14", "17:00", "Bathroom", "No Motion", "Inactive")) | Out-Null
$rows.Add(@("2026-01-30", "17:30:19", "17:00", "Bathroom", "No Motion", "Inactive")) | Out-Null
$rows.Add(@("2026-01-30", "17:30:24", "17:00", "Bathroom", "No Motion", "Inactive")) | Out-Null
$rows.Add(@("2026-01-30", "17:30:29", "17:00", "Bathroom", "No Motion", "Inactive")) | Out-Null
$rows.Add(@("2026-01-30", "17:30:34", "17:00", "Bathroom", "No Motion", "Inactive")) | Out-Null
$rows.Add(@("2026-01-30", "17:30:39", "17:00", "Bathroom", "No Motion", "Inactive")) | Out-Null
$rows.Add(@("2026-01-30", "17:30:44", "17:00", "Bathroom", "No Motion", "Inactive")) | Out-Null
$rows.Add(@("2026-01-30", "17:30:49", "17:00", "Bathroom", "No Motion", "Inactive")) | Out-Null
$rows.Add(@("2026-01-30", "17:30:54", "17:00", "Bathroom", "No Motion", "Inactive")) | Out-Null
$rows.Add(@("2026-01-30", "17:30:59", "17:00", "Bathroom", "No Motion", "Inactive")) | Out-Null
$rows.Add(@("2026-01-30", "17:31:04", "17:00", "Bathroom", "No Motion", "Inactive")) | Out-Null
$rows.Add(@("2026-01-30", "17:31:09", "17:00", "Bathroom", "No Motion", "Inactive")) | Out-Null
$startRow = 278
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $rowNum = $startRow + $i
    $ws.Range("A" + $rowNum + ":F" + $rowNum).NumberFormat = "@"
    for ($j = 0; $j -lt 6; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $r[$j]
    }
}

# ---- Sheet: Humidity (A1:F193 -> A1:F202) ----
$ws = $wb.Worksheets.Item("Humidity")
$rows = New-Object System.Collections.ArrayList
$rows.Add(@("2026-01-30", "17:30:12", "17:00", "Bathroom", "86.2%", "Active")) | Out-Null
$rows.Add(@("2026-01-30", "17:30:13", "17:00", "Bathroom", "87.0%", "Active")) | Out-Null
$rows.Add(@("2026-01-30", "17:30:15", "17:00", "Bathroom", "86.1%", "Active")) | Out-Null
$rows.Add(@("2026-01-30", "17:30:20", "17:00", "Bathroom", "87.0%", "Active")) | Out-Null
$rows.Add(@("2026-01-30", "17:30:40", "17:00", "Bathroom", "87.1%", "Active")) | Out-Null
$rows.Add(@("2026-01-30", "17:30:50", "17:00", "Bathroom", "87.1%", "Active")) | Out-Null
$rows.Add(@("2026-01-30", "17:30:55", "17:00", "Bathroom", "87.1%", "Active")) | Out-Null
$rows.Add(@("2026-01-30", "17:31:00", "17:00", "Bathroom", "87.1%", "Active")) | Out-Null
$rows.Add(@("2026-01-30", "17:31:10", "17:00", "Bathroom", "87.1%", "Active")) | Out-Null
$startRow = 194
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $rowNum = $startRow + $i
    $ws.Range("A" + $rowNum + ":F" + $rowNum).NumberFormat = "@"
    for ($j = 0; $j -lt 6; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $r[$j]
    }
}

# ---- Sheet: mmWave (A1:F52 -> A1:F53) ----
$ws = $wb.Worksheets.Item("mmWave")
$rows = New-Object System.Collections.ArrayList
$rows.Add(@("2026-01-30", "17:30:19", "17:00", "Living Room", "FALL_DETECTED", "EMERGENCY")) | Out-Null
$startRow = 53
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $rowNum = $startRow + $i
    $ws.Range("A" + $rowNum + ":F" + $rowNum).NumberFormat = "@"
    for ($j = 0; $j -lt 6; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $r[$j]
    }
}
